# Rename worksheets and update the label cells that mirror the sheet's
# previous naming ("sw_summary"/"sw_signups" -> "signup_summary_sheet"/"signup_sheet")
# plus the corresponding text values in A1 of each sheet.

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("sw_summary")
$signupsSheet = $wb.Worksheets.Item("sw_signups")

# Update the text values that describe the sheet contents before renaming
# (so lookups by old sheet name still work while we touch the cells).
if ($summarySheet.Range("A1").Value2 -eq "signups_summary_list") {
    $summarySheet.Range("A1").Value = "signup_summary_list"
}

if ($signupsSheet.Range("A1").Value2 -eq "signups_table") {
    $signupsSheet.Range("A1").Value = "signup_table"
}

# Rename the worksheets themselves.
$summarySheet.Name = "signup_summary_sheet"
$signupsSheet.Name = "signup_sheet"
